# Add "Budgeted Labor Units" to the Activities sheet and document it on README.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Activities
# ---------------------------------------------------------------------------
$activities = $wb.Worksheets.Item("Activities")

# Clear the placeholder empty strings left in B2:C2 and B3:C3 (these rows have
# no Activity Name / Activity Status because they are summary/header rows).
$activities.Range("B2:C2").ClearContents()
$activities.Range("B3:C3").ClearContents()

# New column header
$activities.Range("H1").Value = "Budgeted Labor Units"

# New column values
$activities.Range("H2").Value = 1000
$activities.Range("H3").Value = 600
$activities.Range("H4").Value = 300
$activities.Range("H5").Value = 300
$activities.Range("H6").Value = 400

# ---------------------------------------------------------------------------
# Sheet 5: README
# ---------------------------------------------------------------------------
$readme = $wb.Worksheets.Item("README")

# Insert a new bullet point right after "Activities sheet:" (row 3), pushing
# the remaining rows down by one.
$readme.Rows("4:4").Insert()
$readme.Range("A4").Value = "- Include Budgeted Labor Units (last column) for all activities."

# Pad the sheet with trailing blank rows up to row 31 so the used range grows
# to match the template layout, without introducing any new cell styles.
for ($r = 17; $r -le 31; $r++) {
    $readme.Cells.Item($r, 1).Font.Bold = $false
}
